$wb = $excel.ActiveWorkbook

# The "Spain" sheet is the last sheet / template for country-specific test data.
# Duplicate it to create the new "Turkey" sheet (this preserves styles, merged
# cells, column widths, etc. exactly like the source sheet).
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)

# The freshly copied sheet is placed right after "Spain" and is named "Spain (2)".
$turkey = $wb.Worksheets.Item("Spain (2)")
$turkey.Name = "Turkey"

# Fill in the Turkey-specific test data (order matters for shared-string index
# allocation: part number first, then market name).
$turkey.Range("B4").Value = "NGC-3191/T3311 "
$turkey.Range("B2").Value = "Turkey Market"

# Column D was resized on the new sheet (distinct from Spain's column D width).
$turkey.Columns.Item(4).ColumnWidth = 18

# Restore Spain's selection to the full used range (it is no longer the active
# sheet/tab once Turkey is added).
$spain.Range("A1:D14").Select() | Out-Null

# Turkey becomes the newly active sheet/tab, with cell F12 selected.
$turkey.Range("F12").Select() | Out-Null
